# Auto-generated edit script applying market-data refresh values
# per the commit diff (per-sheet H/I/J/K/L/M/N numeric cell updates).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1744.8462
$ws.Range("J17").Value = 1744.8462
$ws.Range("L17").Value = 5234.5386
$ws.Range("N17").Value = -5570.5386

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 3261
$ws.Range("J38").Value = 8999.5
$ws.Range("L38").Value = 26998.5
$ws.Range("N38").Value = -27742.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 940.6842
$ws.Range("I41").Value = 977.0769
$ws.Range("K41").Value = 977.0769
$ws.Range("M41").Value = -537.0769

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 7839.92
$ws.Range("I51").Value = 27650
$ws.Range("K51").Value = 27650
$ws.Range("M51").Value = -27166

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 1949.75
$ws.Range("I53").Value = 1750
$ws.Range("J53").Value = 2149.5
$ws.Range("K53").Value = 1750
$ws.Range("L53").Value = 2149.5
$ws.Range("M53").Value = -1113
$ws.Range("N53").Value = -3423.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 6091.35
$ws.Range("I64").Value = 5408.4443
$ws.Range("K64").Value = 5408.4443
$ws.Range("M64").Value = -5160.4443

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 6091.35
$ws.Range("I67").Value = 5408.4443
$ws.Range("K67").Value = 5408.4443
$ws.Range("M67").Value = -4550.4443

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 12499
$ws.Range("I76").Value = 12499
$ws.Range("K76").Value = 12499
$ws.Range("M76").Value = -12184

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 12499
$ws.Range("I79").Value = 12499
$ws.Range("K79").Value = 12499
$ws.Range("M79").Value = -11407

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 4001.48
$ws.Range("I106").Value = 4172.25
$ws.Range("K106").Value = 4172.25
$ws.Range("M106").Value = -3541.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3444.476
$ws.Range("I122").Value = 3166.7
$ws.Range("J122").Value = 9000
$ws.Range("K122").Value = 9500.099999999999
$ws.Range("L122").Value = 27000
$ws.Range("M122").Value = -7050.099999999999
$ws.Range("N122").Value = -31900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 109900
$ws.Range("J57").Value = 109900
$ws.Range("L57").Value = 109900
$ws.Range("N57").Value = -111340

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 326278.6
$ws.Range("I86").Value = 381450.88
$ws.Range("K86").Value = 381450.88
$ws.Range("M86").Value = -380327.88

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 326278.6
$ws.Range("I89").Value = 381450.88
$ws.Range("K89").Value = 1907254.4
$ws.Range("M89").Value = -1901638.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3574205.8
$ws.Range("I134").Value = 2545.7144
$ws.Range("J134").Value = 14289185
$ws.Range("K134").Value = 7637.1432
$ws.Range("L134").Value = 42867555
$ws.Range("M134").Value = -5102.1432
$ws.Range("N134").Value = -42872625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 109666.664
$ws.Range("J135").Value = 109666.664
$ws.Range("L135").Value = 109666.664
$ws.Range("N135").Value = -119806.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H136").Value = 109900
$ws.Range("J136").Value = 109900
$ws.Range("L136").Value = 109900
$ws.Range("N136").Value = -120100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H141").Value = 191124.12
$ws.Range("J141").Value = 192714
$ws.Range("L141").Value = 192714
$ws.Range("N141").Value = -203074

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 35717090
$ws.Range("I31").Value = 41668900
$ws.Range("K31").Value = 41668900
$ws.Range("M31").Value = -41668605

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 35717090
$ws.Range("I34").Value = 41668900
$ws.Range("K34").Value = 41668900
$ws.Range("M34").Value = -41668698

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1547.3077
$ws.Range("I105").Value = 1191.2727
$ws.Range("J105").Value = 3505.5
$ws.Range("K105").Value = 1191.2727
$ws.Range("L105").Value = 3505.5
$ws.Range("M105").Value = 555.7273
$ws.Range("N105").Value = -6999.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1419.8096
$ws.Range("I107").Value = 554.75
$ws.Range("J107").Value = 2573.2222
$ws.Range("K107").Value = 554.75
$ws.Range("L107").Value = 2573.2222
$ws.Range("M107").Value = 1365.25
$ws.Range("N107").Value = -6413.2222

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2438.8386
$ws.Range("I132").Value = 2218.5652
$ws.Range("J132").Value = 3072.125
$ws.Range("K132").Value = 6655.6956
$ws.Range("L132").Value = 9216.375
$ws.Range("M132").Value = -4125.6956
$ws.Range("N132").Value = -14276.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 166677330
$ws.Range("I12").Value = 500004500
$ws.Range("K12").Value = 1500013500
$ws.Range("M12").Value = -1500013327

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 11591.714
$ws.Range("I76").Value = 9956
$ws.Range("K76").Value = 29868
$ws.Range("M76").Value = -29485

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H79").Value = 11591.714
$ws.Range("I79").Value = 9956
$ws.Range("K79").Value = 29868
$ws.Range("M79").Value = -28542

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1051.1875
$ws.Range("I113").Value = 1103.1666
$ws.Range("J113").Value = 1020
$ws.Range("K113").Value = 3309.4998
$ws.Range("L113").Value = 3060
$ws.Range("M113").Value = -1139.4998
$ws.Range("N113").Value = -7400

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 5195.375
$ws.Range("J121").Value = 7209
$ws.Range("L121").Value = 21627
$ws.Range("N121").Value = -24247

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 758.6667
$ws.Range("I107").Value = 138.5
$ws.Range("K107").Value = 138.5
$ws.Range("M107").Value = 1781.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5255.6113
$ws.Range("I7").Value = 5833.8335
$ws.Range("J7").Value = 4966.5
$ws.Range("K7").Value = 5833.8335
$ws.Range("L7").Value = 4966.5
$ws.Range("M7").Value = -5721.8335
$ws.Range("N7").Value = -5190.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 7789.4
$ws.Range("I22").Value = 13054.556
$ws.Range("J22").Value = 3481.5454
$ws.Range("K22").Value = 13054.556
$ws.Range("L22").Value = 3481.5454
$ws.Range("M22").Value = -12759.556
$ws.Range("N22").Value = -4071.5454

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 7789.4
$ws.Range("I27").Value = 13054.556
$ws.Range("J27").Value = 3481.5454
$ws.Range("K27").Value = 13054.556
$ws.Range("L27").Value = 3481.5454
$ws.Range("M27").Value = -12947.556
$ws.Range("N27").Value = -3695.5454

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 5955581.5
$ws.Range("J68").Value = 5002
$ws.Range("L68").Value = 5002
$ws.Range("N68").Value = -6500

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 5955581.5
$ws.Range("J71").Value = 5002
$ws.Range("L71").Value = 25010
$ws.Range("N71").Value = -32498

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 5255.6113
$ws.Range("I126").Value = 5833.8335
$ws.Range("J126").Value = 4966.5
$ws.Range("K126").Value = 17501.5005
$ws.Range("L126").Value = 14899.5
$ws.Range("M126").Value = -15031.5005
$ws.Range("N126").Value = -19839.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4477.1113
$ws.Range("I126").Value = 4185.143
$ws.Range("K126").Value = 12555.429
$ws.Range("M126").Value = -10085.429
